# Commit: "Added fund docs import call/dist name"
#
# Rename the "Cost Of Investment *" column header (E1) to
# "Face Value For Redemption *", move the current selection to E2,
# and widen column E to fit the new, longer header text.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E1").Value = "Face Value For Redemption *"

# Select E2 (matches the new <selection activeCell="E2" sqref="E2"/>)
$ws.Range("E2").Select()

# Widen column E to accommodate the longer header (was 17.625 -> ~28.9375)
$ws.Columns("E").ColumnWidth = 28.3
